$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("D1").Value = "CredentialName"

# Existing rows 2-4 get a new D column value
$ws.Range("D2").Value = "CredentialUIBANK"
$ws.Range("D3").Value = "CredentialUIBANK"
$ws.Range("D4").Value = "CredentialUIBANK"

# New rows 5-7 (repeat of Primary/Secondary/Test Account, all Exists=FALSE, MinBalance=-1)
$ws.Range("A5").Value = "Primary Account"
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = -1
$ws.Range("D5").Value = "CredentialUIBANKNODATA"

$ws.Range("A6").Value = "Secondary Account"
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = -1
$ws.Range("D6").Value = "CredentialUIBANKNODATA"

$ws.Range("A7").Value = "Test Account"
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = -1
$ws.Range("D7").Value = "CredentialUIBANKNODATA"

# Update selection to match diff
$ws.Range("I16").Select() | Out-Null

# Page setup to match diff (paperSize=9 -> A4, orientation=portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
